# Apply updated bus voltage magnitude (vm_pu) results for the 380 kV case.
# Source bus (B column, the slack/reference bus) setpoint changes from 1.05 to 1.02 p.u.,
# and the downstream bus voltages (C:F and I:N) are refreshed with the recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows correspond to sheet rows 2-25 (results index 0-23).
# Block 1 covers columns B:F.
$data1 = @(
    @(1.02, 1.037247498688274, 1.040038911576404, 1.046185339492755, 1.057464861272936),
    @(1.02, 1.038486082062498, 1.040982579724296, 1.04730713433842, 1.058723065785535),
    @(1.02, 1.039286984818398, 1.041592381204417, 1.04803285964516, 1.059537152704852),
    @(1.02, 1.039623556814952, 1.041848547778179, 1.048337920084305, 1.059879383568444),
    @(1.02, 1.039680061331485, 1.041891547931129, 1.048389139100642, 1.059936845033516),
    @(1.02, 1.039291482607738, 1.041595804875502, 1.048036936011101, 1.05954172565221),
    @(1.02, 1.037666197840369, 1.040357998058212, 1.046564487576098, 1.057890090094765),
    @(1.02, 1.034797959055134, 1.038170531784059, 1.043968613915968, 1.054979180875856),
    @(1.02, 1.032882756908706, 1.03670791527883, 1.042237095696887, 1.053038097385777),
    @(1.02, 1.032052693821785, 1.036073547261255, 1.041487082368079, 1.052197443715698),
    @(1.02, 1.031744252715099, 1.035837756031035, 1.041208453541143, 1.051885162040168),
    @(1.02, 1.031810419805635, 1.0358883412585, 1.041268222238441, 1.051952148715499),
    @(1.02, 1.032027200424902, 1.036054059931655, 1.041464051661182, 1.052171630950163),
    @(1.02, 1.032160750287488, 1.036156143596521, 1.041584703200056, 1.052306857743021),
    @(1.02, 1.032937828997576, 1.036749994029436, 1.042286866057828, 1.053093885376495),
    @(1.02, 1.033425061990233, 1.037122219674288, 1.042727244546185, 1.053587524513196),
    @(1.02, 1.033709182886432, 1.037339231840437, 1.042984085464424, 1.053875441415153),
    @(1.02, 1.033806048341442, 1.037413210340158, 1.043071657508176, 1.053973611266851),
    @(1.02, 1.033372794166961, 1.03708229381968, 1.042679998643563, 1.053534563221308),
    @(1.02, 1.031963367239783, 1.036005264307082, 1.041406385912705, 1.05210699962726),
    @(1.02, 1.031076516058851, 1.035327174073527, 1.040605379700271, 1.051209285023154),
    @(1.02, 1.031546719036404, 1.035686730202064, 1.041030031309727, 1.051685195515361),
    @(1.02, 1.033396411974604, 1.037100334906157, 1.0427013471087, 1.053558494194079),
    @(1.02, 1.035539991676988, 1.038736797603319, 1.044639866466996, 1.055731795228489)
)

# Block 2 covers columns I:N (column G stays 1, column H has no data).
$data2 = @(
    @(1.037404451523114, 1.04235140735136, 1.042822392844397, 1.048951483728923, 1.060199757814831, 1.017972340960573),
    @(1.037712470755133, 1.043232837882837, 1.043576051588511, 1.049884063743063, 1.06127067532175, 1.018271196780265),
    @(1.037909787517211, 1.043802152602022, 1.044062285126043, 1.050486770549699, 1.061963016229517, 1.018464060556817),
    @(1.037992263028824, 1.044041247492964, 1.044266355514999, 1.050739974099959, 1.062253931026522, 1.018545017407291),
    @(1.038006083104972, 1.044081378253121, 1.044300599797646, 1.05078247790804, 1.062302768430091, 1.018558603213676),
    @(1.037910891430572, 1.043805348360539, 1.044065013267762, 1.050490154551282, 1.061966904016119, 1.018465142789306),
    @(1.037508960882382, 1.042649505061393, 1.043077392960677, 1.04926680658232, 1.060561807658061, 1.018073447770037),
    @(1.036785426433533, 1.040604806225259, 1.041326045227452, 1.047105414305955, 1.058081064227437, 1.017379260167118),
    @(1.036292761652366, 1.039236226750822, 1.040150986460655, 1.045660556879343, 1.056423905314185, 1.016913770226659),
    @(1.036076979136745, 1.03864230213143, 1.03964037818607, 1.045033962880994, 1.05570552129141, 1.016711561338837),
    @(1.03599645821302, 1.038421491612389, 1.039450443481675, 1.044801071497493, 1.055438555436301, 1.016636353887137),
    @(1.036013746966529, 1.038468865323137, 1.039491197457727, 1.044851034151331, 1.055495826260129, 1.016652490588448),
    @(1.036070330798237, 1.038624053966518, 1.039624683683416, 1.04501471500658, 1.055683456395455, 1.016705346671517),
    @(1.03610514495316, 1.038719644185562, 1.03970689285399, 1.04511554471254, 1.055799044741807, 1.016737900035513),
    @(1.036307030627144, 1.039275615571597, 1.040184835787111, 1.04570212143212, 1.056471564549637, 1.016927176445686),
    @(1.036433010236212, 1.039624006704417, 1.040484153694303, 1.04606980671719, 1.056893196383852, 1.017045730438469),
    @(1.036506255153428, 1.03982709007227, 1.040658567195293, 1.046284178722657, 1.05713904782586, 1.017114818431177),
    @(1.036531189625843, 1.039896314724269, 1.040718008299637, 1.046357258391848, 1.057222863438641, 1.017138365060792),
    @(1.036419518321889, 1.039586640814569, 1.0404520576813, 1.046030367159856, 1.056847967496077, 1.017033017191158),
    @(1.036053678481757, 1.038578360329558, 1.039585382840138, 1.04496651913437, 1.055628207495543, 1.01668978458453),
    @(1.035821521145516, 1.037943253480697, 1.039038895125065, 1.044296787175738, 1.054860565549187, 1.016473412928247),
    @(1.035944795190389, 1.038280046366456, 1.039328748423068, 1.044651905877747, 1.055267577095022, 1.01658816965336),
    @(1.036425615469087, 1.039603525245078, 1.04046656102434, 1.046048188484278, 1.056868404728545, 1.01703876195365),
    @(1.036974291841828, 1.041134362696547, 1.041780125879739, 1.04766487127593, 1.058722974438767, 1.017559197785714)
)

$rowCount = $data1.Count
$colCount1 = $data1[0].Count
$colCount2 = $data2[0].Count

$block1 = New-Object 'object[,]' $rowCount,$colCount1
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount1; $c++) {
        $block1[$r,$c] = $data1[$r][$c]
    }
}
$ws.Range("B2:F25").Value = $block1

$block2 = New-Object 'object[,]' $rowCount,$colCount2
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount2; $c++) {
        $block2[$r,$c] = $data2[$r][$c]
    }
}
$ws.Range("I2:N25").Value = $block2
